$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 741, shifting existing rows 741:782 down to 742:783
$ws.Rows.Item(741).Insert()

# Populate the newly inserted row with the new data entry.
# Force text format on column A so the date-like string "2026/01/30" is
# stored as a literal string instead of being auto-converted to a date serial,
# then clear the formatting override so the cell keeps the plain, unstyled
# look shared by every other data row.
$ws.Cells.Item(741, 1).NumberFormat = "@"
$ws.Cells.Item(741, 1).Value = "2026/01/30"
$ws.Cells.Item(741, 1).ClearFormats()
$ws.Cells.Item(741, 2).Value = "金"
$ws.Cells.Item(741, 3).Value = 20
$ws.Cells.Item(741, 4).Value = 201
